$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("HaUI", 181.5,  "2025-06-04 14:46:15", "/mnt/d/project/DATN/DATN_PhamDangDong/resource/data/data_result/image_20250604_144608.jpg"),
    @("HaUI", 162.86, "2025-06-05 00:15:30", "/mnt/d/project/DATN/DATN_PhamDangDong/resource/data/data_result/image_20250605_001523.jpg"),
    @("HaUI", 161.76, "2025-06-05 00:51:04", "/mnt/d/project/DATN/DATN_PhamDangDong/resource/data/data_result/image_20250605_005055.jpg"),
    @("HaUI", 175,    "2025-06-05 01:13:52", "/mnt/d/project/DATN/DATN_PhamDangDong/resource/data/data_result/image_20250605_011345.jpg")
)

$startRow = 59
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
